$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: repeat "Start of ramp" / "End of ramp" pairs in F:I
$ws.Range("F1").Value = "Start of ramp"
$ws.Range("G1").Value = "End of ramp"
$ws.Range("H1").Value = "Start of ramp"
$ws.Range("I1").Value = "End of ramp"

# Row 3 - update existing file + values, and extend with new columns
$ws.Range("A3").Value = "/home/daniel/Spike Data/Matlab files/exp 30 slow baseline.mat"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 14230
$ws.Range("E3").Value = 19590
$ws.Range("F3").Value = 124200
$ws.Range("G3").Value = 129200
$ws.Range("H3").Value = 219700
$ws.Range("I3").Value = 225700

# Row 4 - new data row
$ws.Range("A4").Value = "/home/daniel/Spike Data/Matlab files/exp 37 unit 1.mat"
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 7273
$ws.Range("E4").Value = 17560
$ws.Range("F4").Value = 342300
$ws.Range("G4").Value = 350700
$ws.Range("H4").Value = 457100
$ws.Range("I4").Value = 467900

# Row 5 - new data row
$ws.Range("A5").Value = "/home/daniel/Spike Data/Matlab files/exp 43 unit 1.mat"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 3241
$ws.Range("E5").Value = 9703
$ws.Range("F5").Value = 125700
$ws.Range("G5").Value = 135000

# Row 6 - new data row
$ws.Range("A6").Value = "/home/daniel/Spike Data/Matlab files/Exp 41 - two neuronal units slow ramps.mat"
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 13720
$ws.Range("E6").Value = 21230
$ws.Range("F6").Value = 103100
$ws.Range("G6").Value = 110100
$ws.Range("H6").Value = 307600
$ws.Range("I6").Value = 313300

# Update view state: wider tab ratio, scrolled-right view (topLeftCell -> C1),
# and move the selection to I9, matching the author's final cursor position.
$excel.ActiveWindow.TabRatio = 172
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("I9").Select()
